# Atualiza notas dos alunos
# Fill in grade column C2:C7 with the missing student scores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 2.5
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 1.25
$ws.Range("C6").Value = 1.25
$ws.Range("C7").Value = 2.5

# Leave the selection where the author left it after entering the data.
$ws.Range("C8").Select()
